$d = $word.ActiveDocument

function Replace-InRange([int]$startPos, [int]$endPos, [string]$findText, [string]$replaceText) {
    $r = $d.Range($startPos, $endPos)
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# ------------------------------------------------------------------
# Paragraph: "AlgoBay: Utilizo la clase AlgoBay como interfaz..." -> "AlGoOh: Utilizo la clase como interfaz..."
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(48)
$p1Start = $p1.Range.Start
$p1End = $p1.Range.End

Replace-InRange $p1Start $p1End `
    " Utilizo la clase AlgoBay como interfaz entre el cliente y el sistema de compras, pensándola" `
    " Utilizo la clase como interfaz entre el usuario y el modelo, pensándola"

$p1 = $d.Paragraphs.Item(48)
Replace-InRange $p1Start $p1.Range.End `
    " con el usuario (sin tener en cuenta ninguna interfaz gráfica o menú con facilidades)." `
    " (sin tener en cuenta ninguna interfaz gráfica o menú con facilidades) que permitirá al jugador realizar acciones sobre el juego cuando corresponda."

# Rename the italic header "AlgoBay" -> "AlGoOh" (only the title occurrence remains at this point)
$p1 = $d.Paragraphs.Item(48)
Replace-InRange $p1Start $p1.Range.End "AlgoBay" "AlGoOh"

# ------------------------------------------------------------------
# Paragraph: "Compra: Es la clase que se encarga de registrar y..." -> "Fase: Es la clase que se encarga de permitir al Jugador..."
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(50)
$p2Start = $p2.Range.Start

Replace-InRange $p2Start $p2.Range.End "Compra" "Fase"

$p2 = $d.Paragraphs.Item(50)
Replace-InRange $p2Start $p2.Range.End `
    "registrar y contabilizar una venta de uno o varios productos procesando su total según haya agregado envío, garantía o cupón." `
    "permitir al Jugador realizar acciones según sea el momento."

# ------------------------------------------------------------------
# Paragraph: "ServicioExtra: Es la interfaz que añade..." -> "Jugador: Representa a uno de los jugadores..."
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(52)
$p3Start = $p3.Range.Start

Replace-InRange $p3Start $p3.Range.End "ServicioExtra" "Jugador"

$p3 = $d.Paragraphs.Item(52)
Replace-InRange $p3Start $p3.Range.End `
    "Es la interfaz que añade el servicio total a la Compra. La implementan las clases Envío y Garantía (Cupón también; detallo más adelante) y agrega valor al subtotal antes de calcular el precio final." `
    "Representa a uno de los jugadores que interactuará con el programa permitiéndole realizar cambios en el juego, y tendrá su tablero con zonas y su mano con cartas. Conociendo a su oponente, realizará ataques en su tablero."

# ------------------------------------------------------------------
# Paragraph: "Envio: Se añade como Servicio a la Compra..." -> "Monstruo: Es una clase abstracta..."
# ------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(54)
$p4Start = $p4.Range.Start

Replace-InRange $p4Start $p4.Range.End "Envio" "Monstruo"

$p4 = $d.Paragraphs.Item(54)
Replace-InRange $p4Start $p4.Range.End `
    "Se añade como Servicio a la Compra correspondiente según tenga envío con el fin de sumarse al precio total, habiendo añadido Garantía (si corresponde) previamente. No añade costo extra en caso de superar el monto de `$5000." `
    "Es una clase abstracta que define algunos métodos comunes a todos los monstruos y permite a sus subclases sobrecargarlos para realizar las acciones de manera distinta según corresponda a su comportamiento."

Write-Output "Done"
Write-Output "P48: $($d.Paragraphs.Item(48).Range.Text)"
Write-Output "P50: $($d.Paragraphs.Item(50).Range.Text)"
Write-Output "P52: $($d.Paragraphs.Item(52).Range.Text)"
Write-Output "P54: $($d.Paragraphs.Item(54).Range.Text)"
